# Weekly update: insert a new price record as the new row 17,
# pushing the existing rows 17-51 down to 18-52.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44526
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112012
$ws.Range("G17").Value = "Espinaca"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8500
$ws.Range("M17").Value = 8250
$ws.Range("N17").Value = "$/cuna 10 kilos"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 825
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = "Hortaliza"
